$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before the current row 612, shifting rows 612:638 down to 614:640.
$ws.Range("A612:R613").Insert()

# First new row (612): fresh weekly data point, date 2023-05-29 (serial 45075), "Primera" quality.
$ws.Cells.Item(612, 1).Value = 7
$ws.Cells.Item(612, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(612, 3).Value = "Ñuble"
$ws.Cells.Item(612, 4).Value = 45075
$ws.Cells.Item(612, 5).Value = 16
$ws.Cells.Item(612, 6).Value = 100114014
$ws.Cells.Item(612, 7).Value = "Betarraga"
$ws.Cells.Item(612, 8).Value = "Sin especificar"
$ws.Cells.Item(612, 9).Value = "Primera"
$ws.Cells.Item(612, 10).Value = 100
$ws.Cells.Item(612, 11).Value = 900
$ws.Cells.Item(612, 12).Value = 900
$ws.Cells.Item(612, 13).Value = 900
$ws.Cells.Item(612, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(612, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(612, 16).Value = 180
$ws.Cells.Item(612, 17).Value = 5
$ws.Cells.Item(612, 18).Value = "Hortaliza"

# Second new row (613): same date, "Segunda" quality.
$ws.Cells.Item(613, 1).Value = 7
$ws.Cells.Item(613, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(613, 3).Value = "Ñuble"
$ws.Cells.Item(613, 4).Value = 45075
$ws.Cells.Item(613, 5).Value = 16
$ws.Cells.Item(613, 6).Value = 100114014
$ws.Cells.Item(613, 7).Value = "Betarraga"
$ws.Cells.Item(613, 8).Value = "Sin especificar"
$ws.Cells.Item(613, 9).Value = "Segunda"
$ws.Cells.Item(613, 10).Value = 100
$ws.Cells.Item(613, 11).Value = 600
$ws.Cells.Item(613, 12).Value = 600
$ws.Cells.Item(613, 13).Value = 600
$ws.Cells.Item(613, 14).Value = "`$/paquete 5 unidades"
$ws.Cells.Item(613, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(613, 16).Value = 120
$ws.Cells.Item(613, 17).Value = 5
$ws.Cells.Item(613, 18).Value = "Hortaliza"

# Match the date formatting used by the rest of column D.
$ws.Range("D612:D613").NumberFormat = $ws.Range("D611").NumberFormat
